$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.734.58'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '3.322.46'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('D9').Value = '3.318.02'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.181'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.578'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '700.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.78%  '
$ws.Range('D15').Value = '3.859.48'
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '67.711.93'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '3.330.53'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.888'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.16%  '
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.96'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '569.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  +2.33%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '3.701.66'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.50%  '
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('E41').Value = '  +2.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.51%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0667'
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.333'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '131.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.79%  '
